$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 10-32 (old granular rows no longer needed)
$ws.Rows("10:32").Delete()

# Update rows 2-9 with combined tuple-style text
$ws.Range("A2").Value = "('Forest', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('Genju of the Spires', ['{R}', 'Enchantment — Aura', 'Enchant Mountain', '{2}: Enchanted Mountain becomes a 6/1 red Spirit creature until end of turn. It’s still a land.', 'When enchanted Mountain is put into a graveyard, you may return Genju of the Spires from your graveyard to your hand.'])"
$ws.Range("A4").Value = "('Island', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A5").Value = "('Mountain', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A6").Value = "('Okina Nightwatch', ['{4}{G}', 'Creature — Human Monk', 'As long as you have more cards in hand than each opponent, Okina Nightwatch gets +3/+3.', '4/3'])"
$ws.Range("A7").Value = "('Plains', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A8").Value = "('Skyknight Legionnaire', ['{1}{R}{W}', 'Creature — Human Knight', 'Flying, haste', '2/2'])"
$ws.Range("A9").Value = "('Swamp', ['Basic Land — Swamp', '({T}: Add {B}.)'])"
